# Author edit: "Add files via upload" — the uploaded workbook shows the
# "STOCK" column (column C) of Tabla1 bumped to 1 for every data row
# (rows 2-123). A handful of rows were already 1 in the prior version
# (51-54, 75-80, 118-119), so this simply normalises the whole column.
#
# The workbook's last on-screen action was selecting that same C2:C123
# range after editing it (and scrolling so row 107 was in view — that
# scroll position is a pure view/cosmetic detail with no COM-settable
# equivalent here beyond the selection itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize STOCK (column C) for every data row of Tabla1 to 1.
$ws.Range("C2:C123").Value = 1

# Leave the same range selected, matching the saved view state.
$ws.Range("C2:C123").Select()
